# Update introduction in Project Report.
#
# The original document has a single, empty paragraph. We turn it into the
# document title (styled with the built-in "Titel"/Title paragraph style),
# add a blank spacer paragraph, and add the start of the introduction text
# as a new paragraph below it.

$d = $word.ActiveDocument

# The first (and, at this point, only) paragraph in the document.
$titlePara = $d.Paragraphs(1).Range

# Replace its content with the title, styled as "Titel", including the
# proofing-error markers Word leaves around the misspelled "Ferbedienung"
# (missing "n" in "Fernbedienung") exactly as authored.
$titleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
            "<w:pPr><w:pStyle w:val='Titel'/></w:pPr>" +
            "<w:r><w:t xml:space='preserve'>Arduino mit </w:t></w:r>" +
            "<w:proofErr w:type='spellStart'/>" +
            "<w:r><w:t>Ferbedienung</w:t></w:r>" +
            "<w:proofErr w:type='spellEnd'/>" +
            "<w:r><w:t xml:space='preserve'> Steuern &amp; Hinderniserkennung</w:t></w:r>" +
            "</w:p>" +
            "<w:p/>" +
            "<w:p>" +
            "<w:r><w:t xml:space='preserve'>Während ein wenig über </w:t></w:r>" +
            "<w:r><w:t xml:space='preserve'>einem Semester in Automatisation und Robotik habe ich einen Arduino mittels Fernsteuerung steuern können. Nachdem ich ausversehen zu oft gegen </w:t></w:r>" +
            "</w:p>"

$titlePara.InsertXML($titleXml)
